# Rename "Planilha2" to "empresas"
$wb = $excel.ActiveWorkbook
$wsEmpresas = $wb.Worksheets.Item("Planilha2")
$wsEmpresas.Name = "empresas"

$wsAlunos = $wb.Worksheets.Item("alunos")

# Update selection on "alunos" sheet (no longer the active/selected tab, selection stays at G13)
$wsAlunos.Range("G13").Select()

# Make "empresas" the active sheet/tab, and set its selection to A3
$wsEmpresas.Activate()
$wsEmpresas.Range("A3").Select()
